# Applies the "modified and documented tests" edit to the last table of
# TEST DOCUMENTATION.docx:
#   - Table 4, row 14, col 4: "N/A"   -> "F" (curly quotes)
#   - Table 4, row 14, col 5: "true"  -> "true for all vertex"
#   - Table 4, row 15, col 4: "N/A"   -> "F" (curly quotes)
#   - Table 4, row 15, col 5: "false" -> "false for  vertex 3 and 4"
#
# The "true"/"false" cells additionally need to end up split across three
# runs (matching how Word naturally fragments a run when the text is edited
# in place and then appended to), so those are built up with a couple of
# no-op Bold toggles that force a run split without altering formatting.

$d = $word.ActiveDocument
$t = $d.Tables(4)

$LDQ = [char]0x201C   # “
$RDQ = [char]0x201D   # ”

# --- Row 14, Col 4: N/A -> "F" -----------------------------------------
$cell = $t.Cell(14, 4)
$cell.Range.Find.Execute("N/A", $true, $false, $false, $false, $false, $true, 1, $false, "$LDQ" + "F" + "$RDQ", 2) | Out-Null

# --- Row 14, Col 5: true -> true for all vertex (split into 3 runs) ----
$cell = $t.Cell(14, 5)
$rng = $cell.Range
$textEnd = $rng.End - 1
$ins = $d.Range($textEnd, $textEnd)
$ins.InsertAfter(" for all vertex")

$cell = $t.Cell(14, 5)
$s = $cell.Range.Start
# split "rue" away from the appended suffix
$mid = $d.Range($s + 1, $s + 4)
$mid.Bold = 1
$mid.Bold = 0
# split "t" away from "rue"
$first = $d.Range($s, $s + 1)
$first.Bold = 1
$first.Bold = 0

# --- Row 15, Col 4: N/A -> "F" ------------------------------------------
$cell = $t.Cell(15, 4)
$cell.Range.Find.Execute("N/A", $true, $false, $false, $false, $false, $true, 1, $false, "$LDQ" + "F" + "$RDQ", 2) | Out-Null

# --- Row 15, Col 5: false -> false for  vertex 3 and 4 (split into 3 runs)
$cell = $t.Cell(15, 5)
$rng = $cell.Range
$textEnd = $rng.End - 1
$ins = $d.Range($textEnd, $textEnd)
$ins.InsertAfter(" for  vertex 3 and 4")

$cell = $t.Cell(15, 5)
$s = $cell.Range.Start
# split "alse" away from the appended suffix
$mid = $d.Range($s + 1, $s + 5)
$mid.Bold = 1
$mid.Bold = 0
# split "f" away from "alse"
$first = $d.Range($s, $s + 1)
$first.Bold = 1
$first.Bold = 0

Write-Output "done"
